$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.154.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.907.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8338"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3284"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.92%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07072"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08098"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7661"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.896.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.280"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.165.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("E18").Value = "  -1.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007774"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.160.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.031"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1732"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +23.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.318"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.07"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.107"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.370"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.521"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.06031"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.298"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.084"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.272"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7338"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01935"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.795"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4461"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.948"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8579"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.909"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.581"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.55%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.41%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.831"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.061.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.536"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.82%  "
